$d = $word.ActiveDocument

# --- Hunk 1: remove the _GoBack bookmark that currently sits after
#     "angrist_krueger.dta" (Word will re-create it wherever the last
#     edit below actually lands). ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Hunk 2: split the "Add to the previous 2SLS..." run into four runs,
#     inserting "as instruments " after "...year-of-birth indicators "
#     and before "(keeping ...". ---

$quote1 = [char]8220
$quote2 = [char]8221

$full = $d.Content
$found = $full.Find.Execute(
    "Add to the previous 2SLS specification interactions of the three " +
    "quarter-of-birth indicators with all of the year-of-birth indicators " +
    "(keeping the year-of-birth " + $quote1 + "main effects" + $quote2 +
    " as controls). Report your coefficient and standard error. How do " +
    "these compare with the coefficients and standard errors in part 1 " +
    "and 2? Comment on any other 2SLS diagnostics and how they affect " +
    "how you feel about this estimate of the returns to schooling.  ")

if (-not $found) {
    throw "Could not find the target sentence to split"
}

$start = $full.Start

# Remove the original text completely, then rebuild it as four discrete
# runs (each InsertAfter on a freshly-collapsed Range becomes its own run).
$full.Delete()

$piece1 = "Add to the previous 2SLS specification interactions of the three quarter-of-birth indicators with all of the year-of-birth indicators "
$piece2 = "as instruments "
$piece3 = "(keeping the year-of-birth " + $quote1 + "main effects" + $quote2 + " as controls). Report your "
$piece4 = "coefficient and standard error. How do these compare with the coefficients and standard errors in part 1 and 2? Comment on any other 2SLS diagnostics and how they affect how you feel about this estimate of the returns to schooling.  "

$cursor = $d.Range($start, $start)
$cursor.InsertAfter($piece1)

$pos2 = $start + $piece1.Length
$cursor = $d.Range($pos2, $pos2)
$cursor.InsertAfter($piece2)

$pos3 = $pos2 + $piece2.Length
$cursor = $d.Range($pos3, $pos3)
$cursor.InsertAfter($piece3)

$bmPos = $pos3 + $piece3.Length
$cursor = $d.Range($bmPos, $bmPos)
$cursor.InsertAfter($piece4)

# Place the _GoBack bookmark (collapsed) right before "coefficient..." —
# do this as the very last step, after any Range.Delete() calls, since
# a prior non-collapsed delete can otherwise confuse where an explicit
# Bookmarks.Add("_GoBack", ...) lands.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
